$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (workbook.xml <sheet name=".../> and app.xml title)
$ws.Name = "テスト123456"

# Sheet view should no longer be right-to-left
$ws.DisplayRightToLeft = $false

# Row 2, column A: "202473" -> "20240708", keep it stored as text
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "20240708"

# Remove row 3 entirely (identical duplicate of row 2); remaining rows shift
# up and the sheet dimension shrinks from A1:C3 to A1:C2 automatically.
$ws.Rows(3).Delete()
